$d = $word.ActiveDocument

$replacements = @(
    @("58×42=", "90×23="),
    @("46×89=", "47×32="),
    @("59×99=", "12×35="),
    @("50×52=", "77×60="),
    @("79×69=", "67×96="),
    @("68×35=", "93×27="),
    @("17×67=", "42×46="),
    @("60×14=", "87×46="),
    @("76×56=", "31×18="),
    @("35×45=", "92×57="),
    @("88×30=", "26×74="),
    @("21×57=", "18×68="),
    @("77×30=", "19×57="),
    @("99×49=", "78×66="),
    @("63×65=", "96×77="),
    @("29×13=", "39×76="),
    @("35×74=", "38×63="),
    @("43×23=", "32×22="),
    @("50×35=", "44×43="),
    @("33×99=", "42×42="),
    @("53×64=", "70×52="),
    @("30×62=", "18×69="),
    @("44×26=", "68×86="),
    @("82×80=", "25×85="),
    @("43×37=", "35×88=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
